$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 87
$ws.Range("I2").Value = 207
$ws.Range("J2").Value = 898
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 263
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 201
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 13
$ws.Range("S2").Value = 96
$ws.Range("T2").Value = 175
$ws.Range("U2").Value = 17
$ws.Range("V2").Value = 1455
$ws.Range("X2").Value = 1454
$ws.Range("Y2").Value = 0
$ws.Range("AA2").Value = 7
